# Assessment 2 Files added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove the old rows 5-7 (Amazon "addProductToCart" testcase + EOF marker)
#    so only the header row + two testcase rows + EOF marker remain.
# ---------------------------------------------------------------------------
$ws.Rows("5:7").Delete()

# ---------------------------------------------------------------------------
# 2) Replace the contents of row 2 with the "Create Account" testcase.
#    Order matters: it controls the order new shared strings are created in.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Register"
$ws.Range("C2").Value = "Create Account"
$ws.Range("D2").Value = "1. Hit the URL`n2. Validate the Application Language is Changing after Selecting different language other than English and Revert back to English`n3. Click on Account Icon -> Create Account and Enter Information for Signup`n4. Verify Field Validations for Date of birth Field with Future Date of birth [Date 2 days from Today]`n5. Click Create Account Button with Valid information`n6. Click close in Image Selection Popup`n7. Verify Create Account is Reloaded with Entered information and Error Message To Perform Next Action for signup"
$ws.Range("A2").Value = "LanguageChangeToChina"

# ---------------------------------------------------------------------------
# 3) Replace the contents of row 3 with the "Redirect Links" testcase.
# ---------------------------------------------------------------------------
$ws.Range("D3").Value = "1. Hit the URL`n2. Hover on the 1st level Category [ Shop All]`n3. Verify all the displayed links when hovered on Category are properly getting navigate to its respective screen`n4. Verify all the links displayed under `"Help`" Column in Footer Section of the Page are properly getting navigate to its respective screen"
$ws.Range("B3").Value = "Links"
$ws.Range("C3").Value = "Redirecting to Links"
$ws.Range("A3").Value = "RedirectLinks"
$ws.Range("E3").Value = "Application is open"

# ---------------------------------------------------------------------------
# 4) Row 4 becomes the EOF marker row again (values already exist as "EOF").
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "EOF"
$ws.Range("B4").Value = "EOF"
$ws.Range("C4").Value = "EOF"
$ws.Range("D4").Value = "EOF"
$ws.Range("E4").Value = "EOF"

# ---------------------------------------------------------------------------
# 5) Formatting: row 2 -> module-name cell gets a shaded/Consolas look, the
#    summary/testcase-name cells are centred, and the steps cell is centred
#    + wrapped.
# ---------------------------------------------------------------------------
$ws.Range("A2").Font.Name = "Consolas"
$ws.Range("A2").Font.Size = 10
$ws.Range("A2").Interior.Pattern = 1
$ws.Range("A2").Interior.ThemeColor = 2
$ws.Range("A2").VerticalAlignment = -4108

$ws.Range("B2:C2").HorizontalAlignment = -4108
$ws.Range("B2:C2").VerticalAlignment = -4108

$ws.Range("D2").HorizontalAlignment = -4108
$ws.Range("D2").VerticalAlignment = -4108
$ws.Range("D2").WrapText = $true

# D4 drops the word-wrap formatting that the old row had.
$ws.Range("D4").WrapText = $false

Write-Output "done"
